$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is already blank (no cells); clearing it drops the empty <row>
# element from the saved XML without renumbering any of the rows below it.
$ws.Rows.Item(2).ClearContents()

function Set-TextValue($cell, $text) {
    # Plain cell.Value assignment auto-recognizes full "yyyy-mm-dd" style
    # strings as real dates (and stamps a number-format style on the cell).
    # Going through a formula first and then collapsing it to its literal
    # value via Copy / Paste-Values keeps the cell a literal text string
    # (no formula, no date serial, no number-format style) -- matching how
    # the rest of this sheet stores every value.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$rows = @(
    @("2024-05-22", "09:51:09", "Fallo tornillo", "-", "-", "-", "-", "09:51:13", "0:00:04"),
    @("2024-05-22", "09:58:36", "Fallo fijador tapa", "-", "-", "-", "-"),
    @("2024-05-22", "10:04:21", "-", "-", "-", "Robot no coloca bien ferrita", "-")
)

$startRow = 152
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        Set-TextValue $cell $values[$c]
    }
}

$excel.CutCopyMode = $false

# Row 154 (the last row) also carries present-but-empty H154/I154 cells
# (repair time/duration not computed yet for that still-open incident).
# Touching & reverting a formatting property materializes the cell in the
# saved sheet without leaving any visible style on it.
$ws.Cells.Item(154, 8).Font.Bold = $true
$ws.Cells.Item(154, 8).Font.Bold = $false
$ws.Cells.Item(154, 9).Font.Bold = $true
$ws.Cells.Item(154, 9).Font.Bold = $false
